$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "This"
$ws.Range("B1").Value = "Is"
$ws.Range("C1").Value = "A"
$ws.Range("D1").Value = "Test"

$ws.Range("D1").Select()
